# Commit: "Added RMSE table for ADP"
#
# Sheet layout recap:
#   Row 1        : column headers (RMSE Pitch Const/Yaw Const/Pitch Square/...)
#   A2:A5 "Simulation" group -> rows 2-5 are LQR P / LQR PI / LQG PI / ADP
#   A6:A9 "USB" group        -> rows 6-9 are LQR P / LQR PI / LQG PI / ADP
#   A10:A13 "Android" group  -> rows 10-13 are LQR P / LQR PI / LQG PI / ADP
#
# Row 9 is the "USB" simulation's "ADP" controller row. It previously had no
# measured data ("NO DATA" placeholders in C9:H9); this change fills in the
# actual measured RMSE values now that the ADP results are available.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 1.3067
$ws.Range("D9").Value = 6.1990999999999996
$ws.Range("E9").Value = 6.5789999999999997
$ws.Range("F9").Value = 21.192299999999999
$ws.Range("G9").Value = 2.1877
$ws.Range("H9").Value = 3.6307
